$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 09:52"

# Row 11 - Rusia: updated totals
$ws.Range("B11").Value = 99399
$ws.Range("C11").Value = 5841
$ws.Range("D11").Value = 10286
$ws.Range("E11").Value = 88141
$ws.Range("F11").Value = 2300
$ws.Range("G11").Value = 105
$ws.Range("H11").Value = 972

# Rows 28/29 - Singapur overtakes Austria in total cases, so they swap positions
# Row 28 becomes Singapur (updated data)
$ws.Range("A28").Value = "Singapur"
$ws.Range("B28").Value = 15641
$ws.Range("C28").Value = 690
$ws.Range("D28").Value = 1128
$ws.Range("E28").Value = 14499
$ws.Range("F28").Value = 20
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 14

# Row 29 becomes Austria (unchanged data, just moved down)
$ws.Range("A29").Value = "Austria"
$ws.Range("B29").Value = 15357
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 12580
$ws.Range("E29").Value = 2208
$ws.Range("F29").Value = 136
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 569

# Row 59 - Moldavia: updated recuperados/criticos/muertes
$ws.Range("E59").Value = 2556
$ws.Range("G59").Value = 4
$ws.Range("H59").Value = 107

# Row 68 - Uzbekistan: updated activos/recuperados
$ws.Range("D68").Value = 994
$ws.Range("E68").Value = 953

# Rows 92/93 - Letonia overtakes Republica de Chipre in total cases, so they swap positions
# Row 92 becomes Letonia (updated data)
$ws.Range("A92").Value = "Letonia"
$ws.Range("B92").Value = 849
$ws.Range("C92").Value = 13
$ws.Range("D92").Value = 348
$ws.Range("E92").Value = 486
$ws.Range("F92").Value = 4
$ws.Range("G92").Value = 2
$ws.Range("H92").Value = 15

# Row 93 becomes Republica de Chipre (unchanged data, just moved down)
$ws.Range("A93").Value = "Republica de Chipre"
$ws.Range("B93").Value = 837
$ws.Range("C93").Value = 0
$ws.Range("D93").Value = 148
$ws.Range("E93").Value = 674
$ws.Range("F93").Value = 15
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 15

# Row 127 - Tanzania: updated totals/new cases/recuperados
$ws.Range("B127").Value = 306
$ws.Range("C127").Value = 7
$ws.Range("E127").Value = 248
